$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 102.833336
$ws.Range("I9").Value = 95.40000000000001
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 95.40000000000001
$ws.Range("L9").Value = 140
$ws.Range("M9").Value = 73.59999999999999
$ws.Range("N9").Value = -478
$ws.Range("H15").Value = 2496.8372
$ws.Range("I15").Value = 2496.8372
$ws.Range("K15").Value = 7490.5116
$ws.Range("M15").Value = -7321.5116
$ws.Range("H38").Value = 5285.8667
$ws.Range("I38").Value = 948.8
$ws.Range("J38").Value = 13960
$ws.Range("K38").Value = 2846.4
$ws.Range("L38").Value = 41880
$ws.Range("M38").Value = -2474.4
$ws.Range("N38").Value = -42624
$ws.Range("H70").Value = 9860.869000000001
$ws.Range("J70").Value = 11279.421
$ws.Range("L70").Value = 33838.263
$ws.Range("N70").Value = -34378.263
$ws.Range("H73").Value = 9860.869000000001
$ws.Range("J73").Value = 11279.421
$ws.Range("L73").Value = 33838.263
$ws.Range("N73").Value = -35710.263
$ws.Range("H80").Value = 483
$ws.Range("J80").Value = 595.7778
$ws.Range("L80").Value = 1787.3334
$ws.Range("N80").Value = -3783.3334
$ws.Range("H82").Value = 1414.0834
$ws.Range("I82").Value = 1519
$ws.Range("J82").Value = 260
$ws.Range("K82").Value = 4557
$ws.Range("L82").Value = 780
$ws.Range("M82").Value = -4151
$ws.Range("N82").Value = -1592
$ws.Range("H83").Value = 483
$ws.Range("J83").Value = 595.7778
$ws.Range("L83").Value = 5362.000199999999
$ws.Range("N83").Value = -15346.0002
$ws.Range("H85").Value = 1414.0834
$ws.Range("I85").Value = 1519
$ws.Range("J85").Value = 260
$ws.Range("K85").Value = 4557
$ws.Range("L85").Value = 780
$ws.Range("M85").Value = -3153
$ws.Range("N85").Value = -3588
$ws.Range("H98").Value = 3011.524
$ws.Range("I98").Value = 2902.389
$ws.Range("K98").Value = 2902.389
$ws.Range("M98").Value = -1404.389
$ws.Range("H111").Value = 3123.4
$ws.Range("I111").Value = 3107.5
$ws.Range("J111").Value = 3187
$ws.Range("K111").Value = 9322.5
$ws.Range("L111").Value = 9561
$ws.Range("M111").Value = -6255.5
$ws.Range("N111").Value = -15695
$ws.Range("H112").Value = 1798.6471
$ws.Range("J112").Value = 1890.7097
$ws.Range("L112").Value = 5672.1291
$ws.Range("N112").Value = -7888.1291
$ws.Range("H122").Value = 3011.524
$ws.Range("I122").Value = 2902.389
$ws.Range("K122").Value = 8707.167000000001
$ws.Range("M122").Value = -6257.167000000001
$ws.Range("H125").Value = 741.5
$ws.Range("I125").Value = 724.75
$ws.Range("J125").Value = 775
$ws.Range("K125").Value = 6522.75
$ws.Range("L125").Value = 6975
$ws.Range("M125").Value = -4062.75
$ws.Range("N125").Value = -11895
$ws.Range("H132").Value = 17859474
$ws.Range("I132").Value = 20002244
$ws.Range("J132").Value = 3044.5
$ws.Range("K132").Value = 60006732
$ws.Range("L132").Value = 9133.5
$ws.Range("M132").Value = -60004202
$ws.Range("N132").Value = -14193.5
$ws.Range("H137").Value = 4147.778
$ws.Range("I137").Value = 4305.476
$ws.Range("J137").Value = 3595.8333
$ws.Range("K137").Value = 12916.428
$ws.Range("L137").Value = 10787.4999
$ws.Range("M137").Value = -10366.428
$ws.Range("N137").Value = -15887.4999
$ws.Range("H138").Value = 2726.7532
$ws.Range("I138").Value = 1396.125
$ws.Range("J138").Value = 4165.2705
$ws.Range("K138").Value = 4188.375
$ws.Range("L138").Value = 12495.8115
$ws.Range("M138").Value = 951.625
$ws.Range("N138").Value = -22775.8115
$ws.Range("H141").Value = 2015.1282
$ws.Range("I141").Value = 1345.2858
$ws.Range("J141").Value = 7876.25
$ws.Range("K141").Value = 4035.8574
$ws.Range("L141").Value = 23628.75
$ws.Range("M141").Value = 1144.1426
$ws.Range("N141").Value = -33988.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2356.6616
$ws.Range("I32").Value = 2127.8906
$ws.Range("K32").Value = 2127.8906
$ws.Range("M32").Value = -1840.8906
$ws.Range("H61").Value = 1907.0333
$ws.Range("I61").Value = 1527.5927
$ws.Range("J61").Value = 5322
$ws.Range("K61").Value = 1527.5927
$ws.Range("L61").Value = 5322
$ws.Range("M61").Value = -1315.5927
$ws.Range("N61").Value = -5746
$ws.Range("H74").Value = 2110.6858
$ws.Range("I74").Value = 1299.7037
$ws.Range("K74").Value = 1299.7037
$ws.Range("M74").Value = -425.7037
$ws.Range("H77").Value = 2110.6858
$ws.Range("I77").Value = 1299.7037
$ws.Range("K77").Value = 6498.5185
$ws.Range("M77").Value = -2130.5185
$ws.Range("H97").Value = 479.33334
$ws.Range("I97").Value = 429.1
$ws.Range("J97").Value = 579.8
$ws.Range("K97").Value = 429.1
$ws.Range("L97").Value = 579.8
$ws.Range("M97").Value = 66.89999999999998
$ws.Range("N97").Value = -1571.8
$ws.Range("H132").Value = 6318.352
$ws.Range("I132").Value = 3934.209
$ws.Range("K132").Value = 11802.627
$ws.Range("M132").Value = -9272.627
$ws.Range("H136").Value = 1907.0333
$ws.Range("I136").Value = 1527.5927
$ws.Range("J136").Value = 5322
$ws.Range("K136").Value = 4582.7781
$ws.Range("L136").Value = 15966
$ws.Range("M136").Value = -2032.7781
$ws.Range("N136").Value = -21066

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 19998
$ws.Range("J34").Value = 19998
$ws.Range("L34").Value = 19998
$ws.Range("N34").Value = -20226
$ws.Range("H105").Value = 1300.5405
$ws.Range("I105").Value = 1205.4138
$ws.Range("J105").Value = 1645.375
$ws.Range("K105").Value = 1205.4138
$ws.Range("L105").Value = 1645.375
$ws.Range("M105").Value = 541.5862
$ws.Range("N105").Value = -5139.375
$ws.Range("H107").Value = 2252.7297
$ws.Range("I107").Value = 2314.4614
$ws.Range("J107").Value = 2106.818
$ws.Range("K107").Value = 2314.4614
$ws.Range("L107").Value = 2106.818
$ws.Range("M107").Value = -394.4614000000001
$ws.Range("N107").Value = -5946.818
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = $null
$ws.Range("H134").Value = 3862.347
$ws.Range("I134").Value = 1832.0857
$ws.Range("J134").Value = 8938
$ws.Range("K134").Value = 5496.257100000001
$ws.Range("L134").Value = 26814
$ws.Range("M134").Value = -2961.257100000001
$ws.Range("N134").Value = -31884

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 524.5
$ws.Range("I6").Value = 524.5
$ws.Range("K6").Value = 524.5
$ws.Range("M6").Value = -411.5
$ws.Range("H10").Value = 526.2857
$ws.Range("I10").Value = 484.33334
$ws.Range("J10").Value = 778
$ws.Range("K10").Value = 484.33334
$ws.Range("L10").Value = 778
$ws.Range("M10").Value = -345.33334
$ws.Range("N10").Value = -1056
$ws.Range("H31").Value = 2462.9285
$ws.Range("I31").Value = 2277.875
$ws.Range("J31").Value = 3573.25
$ws.Range("K31").Value = 2277.875
$ws.Range("L31").Value = 3573.25
$ws.Range("M31").Value = -1982.875
$ws.Range("N31").Value = -4163.25
$ws.Range("H34").Value = 2462.9285
$ws.Range("I34").Value = 2277.875
$ws.Range("J34").Value = 3573.25
$ws.Range("K34").Value = 2277.875
$ws.Range("L34").Value = 3573.25
$ws.Range("M34").Value = -2075.875
$ws.Range("N34").Value = -3977.25
$ws.Range("H58").Value = 2303.8206
$ws.Range("I58").Value = 2312.0571
$ws.Range("K58").Value = 2312.0571
$ws.Range("M58").Value = -2109.0571
$ws.Range("H105").Value = 9071.591
$ws.Range("I105").Value = 10732.111
$ws.Range("K105").Value = 10732.111
$ws.Range("M105").Value = -8985.111000000001
$ws.Range("H107").Value = 940.34375
$ws.Range("I107").Value = 523
$ws.Range("J107").Value = 1130.0454
$ws.Range("K107").Value = 523
$ws.Range("L107").Value = 1130.0454
$ws.Range("M107").Value = 1397
$ws.Range("N107").Value = -4970.0454
$ws.Range("H122").Value = 1159.303
$ws.Range("I122").Value = 1027.8966
$ws.Range("J122").Value = 2112
$ws.Range("K122").Value = 3083.6898
$ws.Range("L122").Value = 6336
$ws.Range("M122").Value = -633.6898000000001
$ws.Range("N122").Value = -11236
$ws.Range("H132").Value = 834.4
$ws.Range("I132").Value = 724.8387
$ws.Range("J132").Value = 1683.5
$ws.Range("K132").Value = 2174.5161
$ws.Range("L132").Value = 5050.5
$ws.Range("M132").Value = 355.4839000000002
$ws.Range("N132").Value = -10110.5
$ws.Range("H134").Value = 1166.1691
$ws.Range("I134").Value = 1163.017
$ws.Range("J134").Value = 1181.6666
$ws.Range("K134").Value = 3489.051
$ws.Range("L134").Value = 3544.9998
$ws.Range("M134").Value = -954.0510000000004
$ws.Range("N134").Value = -8614.9998
$ws.Range("H136").Value = 2303.8206
$ws.Range("I136").Value = 2312.0571
$ws.Range("K136").Value = 6936.1713
$ws.Range("M136").Value = -4386.1713

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3177.6667
$ws.Range("J3").Value = 3000
$ws.Range("L3").Value = 9000
$ws.Range("N3").Value = -9224
$ws.Range("H4").Value = 28048.6
$ws.Range("I4").Value = 3747
$ws.Range("K4").Value = 11241
$ws.Range("M4").Value = -11129
$ws.Range("H12").Value = 74.39130400000001
$ws.Range("I12").Value = 65.75
$ws.Range("J12").Value = 79
$ws.Range("K12").Value = 197.25
$ws.Range("L12").Value = 237
$ws.Range("M12").Value = -24.25
$ws.Range("N12").Value = -583
$ws.Range("H39").Value = 111
$ws.Range("I39").Value = 111
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 333
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -39
$ws.Range("N39").Value = $null
$ws.Range("H92").Value = 242.28572
$ws.Range("J92").Value = 205.5
$ws.Range("L92").Value = 616.5
$ws.Range("N92").Value = -3112.5
$ws.Range("H115").Value = 236681.33
$ws.Range("J115").Value = 4994.5
$ws.Range("L115").Value = 14983.5
$ws.Range("N115").Value = -17333.5
$ws.Range("H131").Value = 1314.4667
$ws.Range("I131").Value = 1057
$ws.Range("J131").Value = 1408.091
$ws.Range("K131").Value = 3171
$ws.Range("L131").Value = 4224.272999999999
$ws.Range("M131").Value = 1869
$ws.Range("N131").Value = -14304.273
$ws.Range("H134").Value = 3541.5
$ws.Range("I134").Value = 2239.5625
$ws.Range("K134").Value = 6718.6875
$ws.Range("M134").Value = -1648.6875
$ws.Range("H137").Value = 2818.1
$ws.Range("I137").Value = 1000
$ws.Range("J137").Value = 3020.111
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 9060.332999999999
$ws.Range("M137").Value = 2100
$ws.Range("N137").Value = -19260.333
$ws.Range("H139").Value = 3009.6428
$ws.Range("I139").Value = 3126.111
$ws.Range("J139").Value = 2800
$ws.Range("K139").Value = 9378.332999999999
$ws.Range("L139").Value = 8400
$ws.Range("M139").Value = -4238.332999999999
$ws.Range("N139").Value = -18680

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 49799
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 49799
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 49799
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -50111
$ws.Range("H58").Value = 20001
$ws.Range("J58").Value = 20001
$ws.Range("L58").Value = 20001
$ws.Range("N58").Value = -20555
$ws.Range("H80").Value = 3529.1428
$ws.Range("I80").Value = 3402.5
$ws.Range("J80").Value = 3579.8
$ws.Range("K80").Value = 3402.5
$ws.Range("L80").Value = 3579.8
$ws.Range("M80").Value = -2404.5
$ws.Range("N80").Value = -5575.8
$ws.Range("H83").Value = 3529.1428
$ws.Range("I83").Value = 3402.5
$ws.Range("J83").Value = 3579.8
$ws.Range("K83").Value = 17012.5
$ws.Range("L83").Value = 17899
$ws.Range("M83").Value = -12020.5
$ws.Range("N83").Value = -27883
$ws.Range("H102").Value = 2594.4443
$ws.Range("I102").Value = 2550.0417
$ws.Range("J102").Value = 2949.6667
$ws.Range("K102").Value = 2550.0417
$ws.Range("L102").Value = 2949.6667
$ws.Range("M102").Value = -928.0417000000002
$ws.Range("N102").Value = -6193.6667
$ws.Range("H132").Value = 825.6667
$ws.Range("I132").Value = 900.2973
$ws.Range("K132").Value = 2700.8919
$ws.Range("M132").Value = -170.8918999999996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3822.4443
$ws.Range("I7").Value = 3050.25
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 3050.25
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -2938.25
$ws.Range("N7").Value = -10224
$ws.Range("H21").Value = 8251.5
$ws.Range("I21").Value = 2985
$ws.Range("J21").Value = 10007
$ws.Range("K21").Value = 2985
$ws.Range("L21").Value = 10007
$ws.Range("M21").Value = -2811
$ws.Range("N21").Value = -10355
$ws.Range("H44").Value = 33333
$ws.Range("J44").Value = 33333
$ws.Range("L44").Value = 33333
$ws.Range("N44").Value = -34245
$ws.Range("H46").Value = 799.6
$ws.Range("I46").Value = 812
$ws.Range("J46").Value = 750
$ws.Range("K46").Value = 812
$ws.Range("L46").Value = 750
$ws.Range("M46").Value = -624
$ws.Range("N46").Value = -1126
$ws.Range("H82").Value = 43480310
$ws.Range("I82").Value = 76923920
$ws.Range("J82").Value = 3623.8
$ws.Range("K82").Value = 76923920
$ws.Range("L82").Value = 3623.8
$ws.Range("M82").Value = -76923559
$ws.Range("N82").Value = -4345.8
$ws.Range("H85").Value = 43480310
$ws.Range("I85").Value = 76923920
$ws.Range("J85").Value = 3623.8
$ws.Range("K85").Value = 76923920
$ws.Range("L85").Value = 3623.8
$ws.Range("M85").Value = -76922672
$ws.Range("N85").Value = -6119.8
$ws.Range("H126").Value = 3822.4443
$ws.Range("I126").Value = 3050.25
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 9150.75
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -6680.75
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 3109.7407
$ws.Range("I132").Value = 1820.2683
$ws.Range("J132").Value = 7176.5386
$ws.Range("K132").Value = 5460.8049
$ws.Range("L132").Value = 21529.6158
$ws.Range("M132").Value = -2930.8049
$ws.Range("N132").Value = -26589.6158
$ws.Range("H136").Value = 7035.4
$ws.Range("I136").Value = 6124.077
$ws.Range("J136").Value = 8727.857
$ws.Range("K136").Value = 18372.231
$ws.Range("L136").Value = 26183.571
$ws.Range("M136").Value = -15822.231
$ws.Range("N136").Value = -31283.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32304.5
$ws.Range("I54").Value = 22043.834
$ws.Range("J54").Value = 40000
$ws.Range("K54").Value = 22043.834
$ws.Range("L54").Value = 40000
$ws.Range("M54").Value = -21523.834
$ws.Range("N54").Value = -41040
$ws.Range("H100").Value = 1796.3077
$ws.Range("I100").Value = 1533.1111
$ws.Range("J100").Value = 2388.5
$ws.Range("K100").Value = 3066.2222
$ws.Range("L100").Value = 4777
$ws.Range("M100").Value = -2525.2222
$ws.Range("N100").Value = -5859
$ws.Range("H122").Value = 4505.2173
$ws.Range("I122").Value = 4439.7676
$ws.Range("J122").Value = 5443.3335
$ws.Range("K122").Value = 13319.3028
$ws.Range("L122").Value = 16330.0005
$ws.Range("M122").Value = -10869.3028
$ws.Range("N122").Value = -21230.0005
$ws.Range("H126").Value = 1747.2941
$ws.Range("I126").Value = 1739.6364
$ws.Range("K126").Value = 5218.9092
$ws.Range("M126").Value = -2748.9092
$ws.Range("H132").Value = 774.7414
$ws.Range("I132").Value = 518.625
$ws.Range("J132").Value = 2004.1
$ws.Range("K132").Value = 1555.875
$ws.Range("L132").Value = 6012.299999999999
$ws.Range("M132").Value = 974.125
$ws.Range("N132").Value = -11072.3
$ws.Range("H136").Value = 4274.62
$ws.Range("I136").Value = 4647.7905
$ws.Range("J136").Value = 1982.2858
$ws.Range("K136").Value = 13943.3715
$ws.Range("L136").Value = 5946.857400000001
$ws.Range("M136").Value = -11393.3715
$ws.Range("N136").Value = -11046.8574
